$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide column A
$ws.Range("A1").EntireColumn.Hidden = $false

# Select A:B columns
$ws.Range("A1:B1048576").Select()
